# Update the "想去人数" (number of people interested) counts in column F
# for both the "展览" and "全部类型" worksheets, which were regenerated
# by the site build (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map of row -> new value for column F on the affected sheets.
$updates = @{
    2  = 6478
    3  = 186
    5  = 43
    6  = 1945
    10 = 346
    11 = 3
    12 = 5616
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
